$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to store a label row ("種類","債務人",...) in row 1 and a
# single data row (index 13, debtor/creditor/balance/date/reason) in row 2.
# This edit drops the label row, promoting the one land-debt record into
# row 1 so the exported table becomes a single data row (the caller will
# now add "portion"/"total" columns alongside it).
#
# Move C2:G2 -> C1:G1 as *values* (xlPasteValues) so the shared-string
# "1468576" balance keeps its original text type instead of being
# reinterpreted as a number, while the destination keeps its existing
# bordered/centered header style (s=1).
$ws.Range("C2:G2").Copy()
$ws.Range("C1").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# B1 held the "種類" label; it has no replacement value, so blank it out
# while keeping its style.
$ws.Range("B1").Value = ""

# Drop the now-redundant data row (also removes the leading index value
# that lived in A2); nothing else should occupy column A any more.
$ws.Rows.Item(2).Delete()
